# Actualización automática 2025-07-09 12:35:07
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
# HIDALGO HIDALGO PEDRO GUSTAVO / TULCAN NARVAEZ EDITH MARITZA - row 21 (data), row 22 (days-to-fulfil summary)
$wsGrupo.Range("H21").Value = 1593.9     # INODOROS
$wsGrupo.Range("I21").Value = 658.8      # LAVABOS
$wsGrupo.Range("H22").Value = "2 de 20"  # INODOROS days count
$wsGrupo.Range("I22").Value = "2 de 20"  # LAVABOS days count

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F21").Value = 4409.24   # julio - TULCAN NARVAEZ EDITH MARITZA
$wsMensual.Range("F22").Value = 27937.3   # julio - TOTAL

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
# widen the "POR CUMPLIR" column (E, 5th column)
$wsCumpl.Columns.Item(5).ColumnWidth = 23.15

# INODOROS row (7)
$wsCumpl.Range("D7").Value = 2196.9
$wsCumpl.Range("E7").Value = 203.0999999999999
$wsCumpl.Range("F7").Value = 0.915375

# LAVABOS row (8)
$wsCumpl.Range("D8").Value = 687.6
$wsCumpl.Range("E8").Value = -62.60000000000002
$wsCumpl.Range("F8").Value = 1.10016

# TOTAL row (19)
$wsCumpl.Range("D19").Value = 27937.3
$wsCumpl.Range("E19").Value = 37440.69762291769
$wsCumpl.Range("F19").Value = 0.427319603165803
